$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) / Volume(1h) (column E) values for this data refresh.
# Values that look like plain numbers are entered with a leading apostrophe
# so Excel keeps storing them as text, matching the source data's string cells.
$updates = [ordered]@{
    "D2" = '27.439.08'
    "E2" = '  +1.96%  '
    "D3" = '1.836.50'
    "E3" = '  +1.13%  '
    "D4" = "'1.012"
    "E4" = '  +1.02%  '
    "D5" = "'314.46"
    "E5" = '  +1.78%  '
    "E6" = '  +0.93%  '
    "E7" = '  +1.92%  '
    "D8" = "'0.3692"
    "E8" = '  +0.80%  '
    "D9" = "'0.07465"
    "E9" = '  +1.39%  '
    "D10" = "'0.8863"
    "E10" = '  +1.85%  '
    "D11" = "'20.48"
    "E11" = '  +0.97%  '
    "D12" = '1.878.38'
    "E12" = '  +1.42%  '
    "E13" = '  +3.32%  '
    "D14" = "'5.454"
    "E14" = '  +1.43%  '
    "D15" = "'93.06"
    "E15" = '  +1.65%  '
    "D16" = "'6.582"
    "E16" = '  +1.19%  '
    "D18" = "'0.000008822"
    "E18" = '  +1.15%  '
    "D20" = '27.629.75'
    "E20" = '  +2.51%  '
    "E21" = '  +1.02%  '
    "D22" = "'5.322"
    "E22" = '  +0.42%  '
    "D23" = "'10.69"
    "D24" = '2.100.11'
    "E24" = '  +1.03%  '
    "D25" = "'1.911"
    "E25" = '  +0.83%  '
    "D26" = "'152.29"
    "E26" = '  +0.94%  '
    "D27" = "'18.63"
    "E27" = '  +1.74%  '
    "D28" = "'2.140"
    "E28" = '  +0.40%  '
    "D29" = "'5.254"
    "E29" = '  -0.17%  '
    "D30" = "'117.60"
    "E30" = '  +1.89%  '
    "D31" = "'0.08997"
    "E31" = '  +1.18%  '
    "D32" = "'0.7569"
    "E32" = '  +0.07%  '
    "D33" = "'1.177"
    "E33" = '  +1.62%  '
    "D34" = "'4.554"
    "E34" = '  +1.61%  '
    "D35" = "'2.945"
    "E35" = '  +1.21%  '
    "E36" = '  +1.09%  '
    "D37" = "'1.103"
    "E37" = '  +1.73%  '
    "D38" = "'0.05347"
    "E38" = '  +1.18%  '
    "D39" = "'0.01956"
    "E39" = '  +0.40%  '
    "D40" = "'2.980"
    "E40" = '  -0.25%  '
    "E41" = '  +1.46%  '
    "D42" = "'2.413"
    "E42" = '  +4.56%  '
    "D43" = "'0.5334"
    "E43" = '  +0.61%  '
    "E44" = '  +0.29%  '
    "D45" = "'8.515"
    "E45" = '  +0.88%  '
    "D46" = "'0.4924"
    "E46" = '  +1.20%  '
    "D47" = "'10.59"
    "E47" = '  +1.98%  '
    "E48" = '  +1.06%  '
    "D49" = "'104.94"
    "E49" = '  +1.48%  '
    "E50" = '  +1.00%  '
    "D51" = "'0.06307"
    "E51" = '  +0.25%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
